$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "dnasr281@gmail.com, "
$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $txt = $cell.Text

    if ($txt -and $txt.StartsWith($prefix)) {
        $rest = $txt.Substring($prefix.Length)
        $swapped = $rest + ", dnasr281@gmail.com"
        $cell.Value2 = $swapped
        $changed++
    }
}

Write-Host ("Rows changed: " + $changed)
